$d = $word.ActiveDocument

$apos = [char]0x2019

# ---------------------------------------------------------------------------
# 1) "The True value mean ..." paragraph sentence tweak (do this BEFORE the
#    IF formula / bookmark edit below, so the _GoBack bookmark ends up on the
#    later edit, matching natural editing order).
#    "that the this is date" -> "that this is date"
#    "went to an empty cell the typed this IF formula:" ->
#       "went to an empty cell then typed:"
# ---------------------------------------------------------------------------
$oldSentence = "The True value mean that the this is date, false mean that this isn" + $apos + "t a date, went to an empty cell the typed this IF formula:"
$newSentence = "The True value mean that this is date, false mean that this isn" + $apos + "t a date, went to an empty cell then typed:"
$rng = $d.Content
$ok1 = $rng.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)

# ---------------------------------------------------------------------------
# 2) Fix the split "=DATE(...)" formula: join the two runs (which used to be
#    split around a bookmark) back into one continuous formula, dropping the
#    bookmark that used to sit in the middle of it.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$ok2 = $rng2.Find.Execute("=DATE(RIGHT(E6,2)+2000,MID(E6,4,2),LEFT(E6,2))", $true, $false, $false, $false, $false, $true, 1, $false, "=DATE(RIGHT(E6,2)+2000,MID(E6,4,2),LEFT(E6,2))", 2)

# ---------------------------------------------------------------------------
# 3) Update the IF formula text.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$ok3 = $rng3.Find.Execute("=IF(J2= TRUE,E2,I2)", $true, $false, $false, $false, $false, $true, 1, $false, "=IF(ISNUMBER(J2),E2,I2)", 2)

# ---------------------------------------------------------------------------
# 4) Remove the stray leading space on "and pulled it to all cells".
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$ok4 = $rng4.Find.Execute(" and pulled it to all cells", $true, $false, $false, $false, $false, $true, 1, $false, "and pulled it to all cells", 2)

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark to sit right before "and pulled it to all
#    cells" -- this is where Word leaves it after the most recent edit.
# ---------------------------------------------------------------------------
$rng5 = $d.Content
$ok5 = $rng5.Find.Execute("and pulled it to all cells", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5.Collapse(1)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $rng5)

Write-Host "Edits applied:" $ok1 $ok2 $ok3 $ok4 $ok5
